# Apply "updated with tests, new code (lstm)" edits to the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (row 1) renames ---
$ws.Range("C1").Value  = "GDP"
$ws.Range("E1").Value  = "Budget_Previous_Year"
$ws.Range("F1").Value  = "LatinAmerica"
$ws.Range("G1").Value  = "Africa"
$ws.Range("H1").Value  = "Confessional"
$ws.Range("I1").Value  = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Column C numeric value updates (rows 2-50, recomputed GDP figures) ---
$ws.Range("C2").Value  = 2870.311589353206
$ws.Range("C3").Value  = 1909.084588129339
$ws.Range("C4").Value  = 14239.03920301361
$ws.Range("C5").Value  = 4132.902312418774
$ws.Range("C6").Value  = 2100.656463590606
$ws.Range("C7").Value  = 1268.249210347625
$ws.Range("C8").Value  = 1250.795760575873
$ws.Range("C9").Value  = 471.9591970298227
$ws.Range("C10").Value = 1286.515571617672
$ws.Range("C11").Value = 892.5687203369533
$ws.Range("C12").Value = 2812.435974421079
$ws.Range("C13").Value = 2898.942214704482
$ws.Range("C14").Value = 1955.461557360978
$ws.Range("C15").Value = 13825.35808833117
$ws.Range("C16").Value = 4550.453595838572
$ws.Range("C17").Value = 1357.563719132622
$ws.Range("C18").Value = 1317.890706178356
$ws.Range("C19").Value = 863.7612548677739
$ws.Range("C20").Value = 2828.483778716848
$ws.Range("C21").Value = 2024.117324382548
$ws.Range("C22").Value = 4961.234688573883
$ws.Range("C23").Value = 1410.426304742003
$ws.Range("C24").Value = 2860.874335573629
$ws.Range("C25").Value = 909.3123437708064
$ws.Range("C26").Value = 2094.024217383061
$ws.Range("C27").Value = 14735.09353649063
$ws.Range("C28").Value = 5325.160106166602
$ws.Range("C29").Value = 1469.177610078392
$ws.Range("C30").Value = 2379.668184479739
$ws.Range("C31").Value = 1443.492614888721
$ws.Range("C32").Value = 2887.250212489506
$ws.Range("C33").Value = 929.4690557368662
$ws.Range("C34").Value = 2201.396847776877
$ws.Range("C35").Value = 14721.85595470026
$ws.Range("C36").Value = 5710.587873377512
$ws.Range("C37").Value = 1544.619247249133
$ws.Range("C38").Value = 2497.68592515536
$ws.Range("C39").Value = 1505.810948829135
$ws.Range("C40").Value = 1401.753174264641
$ws.Range("C41").Value = 961.3778847738438
$ws.Range("C42").Value = 6103.590270484282
$ws.Range("C43").Value = 3008.669179463094
$ws.Range("C44").Value = 1640.18070024053
$ws.Range("C45").Value = 14025.35756477021
$ws.Range("C46").Value = 1441.783971398429
$ws.Range("C47").Value = 956.659691840205
$ws.Range("C48").Value = 6500.281937297324
$ws.Range("C49").Value = 3012.536723186288
$ws.Range("C50").Value = 1751.664428859304

# Rows 51 and 56 previously held the text placeholder ".." and now
# become numeric 0.
$ws.Range("C51").Value = 0
$ws.Range("C52").Value = 2425.561644739583
$ws.Range("C53").Value = 1469.192636109792
$ws.Range("C54").Value = 6907.962010581965
$ws.Range("C55").Value = 2854.757682901436
$ws.Range("C56").Value = 0
$ws.Range("C57").Value = 1752.531946133768
